$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows below the existing data row (row 16) to make room for
# the additional debtor rows, preserving the later "Novedad"/signature block.
$ws.Rows("17:18").Insert()

# Duplicate the original data row (now still the template row 16) down into
# the two freshly inserted rows so they inherit identical cell styling
# (borders/number formats) to the source row.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# Row 16: new worker NOLBERTO MATEUS RODRIGUEZ, period 2502
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1049018193"
$ws.Range("D16").Value = "NOLBERTO MATEUS RODRIGUEZ"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 60000
$ws.Range("G16").Value = 1500000

# Row 17: existing worker FIDEL ANDRES CARVAJAL HERNANDEZ, period 2504 (unchanged values)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1065292180"
$ws.Range("D17").Value = "FIDEL ANDRES CARVAJAL HERNANDEZ"
$ws.Range("E17").Value = "2504"
$ws.Range("F17").Value = 16000
$ws.Range("G17").Value = 2000000

# Row 18: same worker FIDEL ANDRES CARVAJAL HERNANDEZ, new period 2503
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1065292180"
$ws.Range("D18").Value = "FIDEL ANDRES CARVAJAL HERNANDEZ"
$ws.Range("E18").Value = "2503"
$ws.Range("F18").Value = 2667
$ws.Range("G18").Value = 2000000

# Update the summary totals: total overdue value, worker count, period count
$ws.Range("E11").Value = 78667
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3
